$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134, shifting existing rows 134:192 down to 135:193
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with a copy of the (now shifted) row 135's data,
# except for the Fecha (date) column which gets a new value.
$ws.Cells.Item(134, 1).Value = 3
$ws.Cells.Item(134, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(134, 3).Value = "Coquimbo"
$ws.Cells.Item(134, 4).Value = 44489
$ws.Cells.Item(134, 5).Value = 5
$ws.Cells.Item(134, 6).Value = 100112039
$ws.Cells.Item(134, 7).Value = "Ciboulette"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 160
$ws.Cells.Item(134, 11).Value = 1500
$ws.Cells.Item(134, 12).Value = 1500
$ws.Cells.Item(134, 13).Value = 1500
$ws.Cells.Item(134, 14).Value = "`$/docena de atados"
$ws.Cells.Item(134, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(134, 16).Value = 500
$ws.Cells.Item(134, 17).Value = 3
$ws.Cells.Item(134, 18).Value = "Hortaliza"
